# Error Calculations and Plots
# Applies the missing-data edits described by the commit: two whole data
# rows (RM 232 and SC 92) are removed from the table, and a number of
# individual cells gain/lose a value (simulating newly-imputed /
# newly-missing measurements).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the two rows that disappear from the table ------------------
# Delete the lower-indexed row last so the earlier delete doesn't shift the
# target row out from under us.
# Before the edit:
#   row 26 = "RM 232"
#   row 28 = "SC 92"
$ws.Rows.Item(28).Delete()   # "SC 92"
$ws.Rows.Item(26).Delete()   # "RM 232"

# --- Per-cell value changes (row numbers are POST-deletion) -------------
$ws.Range("F2").ClearContents()
$ws.Range("C3").Value = 11.2
$ws.Range("F3").ClearContents()
$ws.Range("D4").ClearContents()
$ws.Range("F4").ClearContents()
$ws.Range("C5").ClearContents()
$ws.Range("F5").Value = 17.66
$ws.Range("E6").Value = -5.7
$ws.Range("D9").Value = -14.5
$ws.Range("D10").Value = -14.7
$ws.Range("E12").ClearContents()
$ws.Range("F13").ClearContents()
$ws.Range("E14").Value = -5.4
$ws.Range("F15").Value = 16.2
$ws.Range("D17").ClearContents()
$ws.Range("E17").Value = -7.3
$ws.Range("D18").ClearContents()
$ws.Range("E19").Value = -6.5
$ws.Range("E20").ClearContents()
$ws.Range("C21").Value = 12.7
$ws.Range("C23").ClearContents()
$ws.Range("E23").ClearContents()
$ws.Range("F23").ClearContents()
$ws.Range("E25").ClearContents()
$ws.Range("E27").Value = -10
$ws.Range("E28").Value = -5.9
$ws.Range("F31").Value = 17.18
$ws.Range("C32").Value = 10.5
